$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.239.06"
$ws.Range("E2").Value = "  +1.35%  "
# Row 3
$ws.Range("D3").Value = "2.366.15"
$ws.Range("E3").Value = "  +1.34%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("D5").Value = "'548.47"
$ws.Range("E5").Value = "  +1.40%  "
# Row 6
$ws.Range("D6").Value = "'139.78"
$ws.Range("E6").Value = "  +2.67%  "
# Row 7
$ws.Range("E7").Value = "  +0.02%  "
# Row 8
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  -0.40%  "
# Row 9
$ws.Range("D9").Value = "2.369.16"
$ws.Range("E9").Value = "  +1.46%  "
# Row 10
$ws.Range("E10").Value = "  +3.23%  "
# Row 11
$ws.Range("E11").Value = "  +1.82%  "
# Row 12
$ws.Range("E12").Value = "  +1.34%  "
# Row 13
$ws.Range("E13").Value = "  +3.31%  "
# Row 14
$ws.Range("D14").Value = "'25.61"
$ws.Range("E14").Value = "  +4.90%  "
# Row 15
$ws.Range("E15").Value = "  +7.71%  "
# Row 16
$ws.Range("D16").Value = "2.794.49"
$ws.Range("E16").Value = "  +1.40%  "
# Row 17
$ws.Range("D17").Value = "61.137.41"
$ws.Range("E17").Value = "  +1.11%  "
# Row 18
$ws.Range("D18").Value = "2.365.94"
$ws.Range("E18").Value = "  +1.32%  "
# Row 19
$ws.Range("E19").Value = "  +5.13%  "
# Row 20
$ws.Range("E20").Value = "  +2.47%  "
# Row 21
$ws.Range("D21").Value = "'321.20"
$ws.Range("E21").Value = "  +1.51%  "
# Row 22
$ws.Range("D22").Value = "'6.64"
$ws.Range("E22").Value = "  +1.54%  "
# Row 23
$ws.Range("E23").Value = "  +0.11%  "
# Row 24
$ws.Range("D24").Value = "'64.30"
$ws.Range("E24").Value = "  +2.24%  "
# Row 25
$ws.Range("D25").Value = "'1.72"
$ws.Range("E25").Value = "  -6.77%  "
# Row 26
$ws.Range("D26").Value = "'8.88"
$ws.Range("E26").Value = "  +3.62%  "
# Row 28
$ws.Range("D28").Value = "'536.20"
$ws.Range("E28").Value = "  +7.58%  "
# Row 30
$ws.Range("D30").Value = "'8.23"
$ws.Range("E30").Value = "  +3.98%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0905"
$ws.Range("E31").Value = "  +1.95%  "
# Row 32
$ws.Range("E32").Value = "  +0.95%  "
# Row 33
$ws.Range("E33").Value = "  +2.08%  "
# Row 34
$ws.Range("E34").Value = "  +3.24%  "
# Row 35
$ws.Range("E35").Value = "  -0.38%  "
# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.07%  "
# Row 37
$ws.Range("D37").Value = "'5.62"
$ws.Range("E37").Value = "  +8.00%  "
# Row 38
$ws.Range("D38").Value = "'4.68"
$ws.Range("E38").Value = "  +2.10%  "
# Row 39
$ws.Range("D39").Value = "'1.90"
$ws.Range("E39").Value = "  +6.07%  "
# Row 40
$ws.Range("E40").Value = "  +2.17%  "
# Row 41
$ws.Range("E41").Value = "  +1.41%  "
# Row 42
$ws.Range("D42").Value = "'145.76"
$ws.Range("E42").Value = "  +6.04%  "
# Row 43
$ws.Range("E43").Value = "  +0.09%  "
# Row 44
$ws.Range("E44").Value = "  +3.44%  "
# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.23"
$ws.Range("E45").Value = "  +6.62%  "
# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'147.04"
$ws.Range("E46").Value = "  +4.38%  "
# Row 47
$ws.Range("E47").Value = "  +2.15%  "
# Row 49
$ws.Range("D49").Value = "'20.10"
$ws.Range("E49").Value = "  +3.95%  "
# Row 50
$ws.Range("E50").Value = "  +2.58%  "
# Row 51
$ws.Range("D51").Value = "'0.0904"
$ws.Range("E51").Value = "  +0.87%  "
